$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Read the original (row-oriented) data before we start overwriting it ---
# Row layout (before):
#   A1 = "year" (label), B1 = "exchange_rate" (label)
#   A2:A10 = years 2008..2016, B2:B10 = exchange rate values

$years  = @()
$rates  = @()
for ($r = 2; $r -le 10; $r++) {
    $years += $ws.Cells.Item($r, 1).Value2
    $rates += $ws.Cells.Item($r, 2).Value2
}

# --- Remove the old A3:B10 block entirely (collapses dimension / style pool) ---
$ws.Range("A3:B10").Clear()

# --- Rebuild as a wide (pivoted) table: labels in column A, years across row 1, values across row 2 ---
# A1 keeps the "year" label, A2 gets the "exchange_rate" label (was in B1)
$ws.Range("A2").Value2 = $ws.Range("B1").Value2
$ws.Range("B1").Clear()

for ($i = 0; $i -lt $years.Count; $i++) {
    $col = $i + 2   # B = 2, C = 3, ...
    $ws.Cells.Item(1, $col).Value2 = $years[$i]
    $ws.Cells.Item(2, $col).Value2 = $rates[$i]

    $ws.Cells.Item(1, $col).HorizontalAlignment = -4152
    $ws.Cells.Item(1, $col).VerticalAlignment = -4108

    $ws.Cells.Item(2, $col).HorizontalAlignment = -4152
    $ws.Cells.Item(2, $col).VerticalAlignment = -4108
    $ws.Cells.Item(2, $col).NumberFormat = "#,##0.##"
}

# --- Styles for the label cells ---
# A2 ("exchange_rate") matches the old right+vcenter style (same as the year header used to have)
$ws.Range("A2").HorizontalAlignment = -4152
$ws.Range("A2").VerticalAlignment = -4108

# A1 ("year") gets a new right-only (no vertical-center) alignment style
$ws.Range("A1").HorizontalAlignment = -4152
$ws.Range("A1").VerticalAlignment = -4107

# --- Selection matches the post-edit cursor position recorded in the workbook ---
$ws.Range("J3").Select()
